# pelada_sabado_2024_02.xlsx update
# - fix two mis-typed "jorge" entries to "Jorge" (the lowercase shared-string
#   entry disappears on save once nothing refers to it any more, which also
#   shifts every other shared-string index down by one - matching the diff)
# - append 22 new player rows (421-442), two of which introduce brand new
#   shared strings ("Lucas" and "joão")
# - grow the AutoFilter / _FilterDatabase range to the new data extent
# - move the selection to the first empty row below the new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the two "jorge" typos -> "Jorge" -------------------------------
$ws.Range("A266").Value = "Jorge"
$ws.Range("A365").Value = "Jorge"

# --- grow the autofilter / filter-database range to A1:K421 -------------
# (done BEFORE appending the new rows below - otherwise the engine re-derives
#  the autofilter range from the sheet's full used range instead of the
#  literal range we pass in)
try { $ws.AutoFilterMode = $false } catch {}
$ws.Range("A1:K421").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Jogadores!_FilterDatabase") {
        $n.RefersTo = "=Jogadores!`$A`$1:`$K`$421"
    }
}

# --- append the new rows -------------------------------------------------
$newRows = @(
    @{ Row = 421; Name = "Jorge";        Vals = @(2,2,3,1,1,0,0,0,0) },
    @{ Row = 422; Name = "Eduardo";      Vals = @(2,2,3,0,1,0,0,0,0) },
    @{ Row = 423; Name = "Leandrão";     Vals = @(2,2,3,0,1,0,0,0,0) },
    @{ Row = 424; Name = "Adriano";      Vals = @(2,2,3,5,1,0,0,0,0) },
    @{ Row = 425; Name = "Coxinha";      Vals = @(2,2,3,0,1,0,0,0,0) },
    @{ Row = 426; Name = "Corinthiano";  Vals = @(3,2,3,0,1,0,0,0,0) },
    @{ Row = 427; Name = "Digão";        Vals = @(3,2,3,0,1,0,0,0,0) },
    @{ Row = 428; Name = "Du";           Vals = @(3,2,3,1,1,0,0,0,0) },
    @{ Row = 429; Name = "Marcos";       Vals = @(3,2,3,1,1,0,0,0,0) },
    @{ Row = 430; Name = "Fabinho";      Vals = @(3,2,3,4,1,0,0,0,0) },
    @{ Row = 431; Name = "Lucas";        Vals = @(5,2,2,0,1,1,0,0,0) },
    @{ Row = 432; Name = "Juscielio";    Vals = @(5,2,2,2,1,1,0,0,0) },
    @{ Row = 433; Name = "Leandrinho";   Vals = @(5,2,2,3,1,1,0,0,0) },
    @{ Row = 434; Name = "Marcelão";     Vals = @(5,2,2,1,1,1,0,0,0) },
    @{ Row = 435; Name = "Ismael";       Vals = @(5,2,2,1,1,1,0,0,0) },
    @{ Row = 436; Name = "Guinha";       Vals = @(1,4,2,0,1,0,1,0,0) },
    @{ Row = 437; Name = "joão";         Vals = @(1,4,2,1,1,0,1,0,0) },
    @{ Row = 438; Name = "Athos";        Vals = @(1,4,2,0,1,0,1,0,0) },
    @{ Row = 439; Name = "Cabeleira";    Vals = @(1,4,2,3,1,0,1,0,0) },
    @{ Row = 440; Name = "Peixe";        Vals = @(1,4,2,0,1,0,1,0,0) },
    @{ Row = 441; Name = "Matheus";      Vals = @(4,5,5,0,1,1,0,1,15) },
    @{ Row = 442; Name = "Lucian";       Vals = @(5,5,4,0,1,0,1,0,10) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Name
    $vals = $entry.Vals
    # columns C..K are indices 3..11
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, 3 + $i).Value = $vals[$i]
    }
}

# --- move selection to the next empty row after the appended data -------
$ws.Range("A443").Select() | Out-Null
